$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B24 was stored as an inline/text string "2" -> convert to a real numeric 2
$ws.Range("B24").Value = 2

# Append new row 25 with the new annotation data
$ws.Range("A25").Value = "Ruilin"

# B25 must stay textual ("3"), not auto-converted to a number by Excel
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "3"
$ws.Range("B25").Style = "Normal"

$ws.Range("C25").Value = "无"
$ws.Range("D25").Value = "FBK"
$ws.Range("E25").Value = "OTH"
$ws.Range("F25").Value = "91b1b71f-4957-400a-bdb5-bced2ed448de"
$ws.Range("G25").Value = "S1CChZ-CZ_annotated.xlsx"
$ws.Range("H25").Value = "It took us as a few weeks to reply because we took the time to implement as much as possible of the feedback."
